$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell text: replace "(location)" with "(Locations)" in the relevant rows
$ws.Range("A7").Value = "Highest Covid Infected Cases per country(Locations)"
$ws.Range("A8").Value = "Highest Motality Count due to Covid per country(Locations)"
$ws.Range("A12").Value = "Vaccination by running total per country(Locations)"
$ws.Range("A13").Value = "Percentage for Running total of Vaccinated Vs Population per country(Locations)"

# Update the selected cell to A7, matching the saved view state
$ws.Range("A7").Select()
